$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill in the "TIPO" (type) column B for rows that were missing it.
# host_response_time, host_response_rate, host_acceptance_rate, host_verifications -> object
$objectRows = @(29, 30, 31, 38)
foreach ($r in $objectRows) {
    $ws.Cells.Item($r, 2).Value = "object"
}

# review_scores_* rows -> float
$floatRows = @(75, 76, 77, 78, 79, 80, 81)
foreach ($r in $floatRows) {
    $ws.Cells.Item($r, 2).Value = "float"
}

# Apply an AutoFilter over the data table (header row 13 through last row 88)
$ws.Range("A13:E88").AutoFilter() | Out-Null

# AutoFilter normally registers a hidden workbook-level "_FilterDatabase" name
# scoped to the sheet; recreate that bookkeeping explicitly.
$filterName = $ws.Names.Add('_xlnm._FilterDatabase', '=Hoja1!$A$13:$E$88')
$filterName.Visible = $false

# Update the view: scroll so row 4 is the top-left visible row, set zoom, and
# move the active selection to C20.
$ws.Activate()
$window = $excel.ActiveWindow
$window.ScrollRow = 4
$window.Zoom = 81
$ws.Range("C20").Select()
